$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

# Row 2: Exam_Prerequisit_for_iProc_TC_ID_166 / Login with diff venu
$ws.Range("A2").Value = "Exam_Prerequisit_for_iProc_TC_ID_166"
$ws.Range("B2").Value = "@iProctorRegression Verify Elumina Login with diff venu"

# Row 3: Exam_Prerequisit_for_iProc_TC_ID_166 / Registration and adding existing user
$ws.Range("A3").Value = "Exam_Prerequisit_for_iProc_TC_ID_166"
$ws.Range("B3").Value = "@iProctorRegression Verify Elumina Registration and adding existing user"

# Row 4: iProc_TC_ID_166 / Validation of Candidate attends exam within browser screen and not in Full-screen mode
$ws.Range("A4").Value = "iProc_TC_ID_166"
$ws.Range("B4").Value = "@iProctorRegression Validation of Candidate attends exam within browser screen and not in Full-screen mode"
